{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer paragraphs\n// (and the blank paragraph immediately preceding them) that followed the\n// bibliography text, per the site rebuild that dropped the scraped footer.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the bibliography paragraph that should now be immediately followed\n// by the pre-existing trailing blank paragraph / page-break paragraph.\nconst marker = \"Artigos de revistas especializadas e de jornais; Estudos, artigos, not\u00edcias e pesquisas via internet.\";\nlet markerIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === marker) {\n    markerIndex = i;\n    break;\n  }\n}\n\nif (markerIndex !== -1) {\n  const targets = [\n    \"\",\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\",\n  ];\n\n  for (let t = 0; t < targets.length; t++) {\n    const idx = markerIndex + 1 + t;\n    if (idx >= items.length) break;\n    if (items[idx].text.trim() !== targets[t]) break;\n    items[idx].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer paragraphs\n# (and the blank paragraph immediately preceding them) that followed the\n# bibliography text, per the site rebuild that dropped the scraped footer.\n$d = $word.ActiveDocument\n\n$marker = \"Artigos de revistas especializadas e de jornais; Estudos, artigos, not\u00edcias e pesquisas via internet.\"\n\n$markerIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $marker) {\n        $markerIndex = $i\n        break\n    }\n}\n\nif ($markerIndex -ne -1) {\n    $targets = @(\n        \"\",\n        \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n        \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n    )\n\n    # Delete from the last target back to the first so earlier paragraph\n    # indices stay valid while later ones are removed.\n    for ($t = $targets.Count - 1; $t -ge 0; $t--) {\n        $pIdx = $markerIndex + 1 + $t\n        if ($pIdx -gt $d.Paragraphs.Count) { continue }\n        $p = $d.Paragraphs.Item($pIdx)\n        $txt = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($txt -eq $targets[$t]) {\n            $p.Range.Delete()\n        }\n    }\n}\n"}
